$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Austria (row 2): Kapanma 4 Start/End were free-text dates, now stored
# as real date values. H2 also becomes left aligned like the other date
# columns (matches column G's style).
$ws.Range("H2").Value = 44522
$ws.Range("H2").HorizontalAlignment = -4131
$ws.Range("I2").Value = 44542

# --- Finland (row 9): the Kapanma-1 dates that were sitting in B9/C9 really
# belong in D9/E9, so copy them across first (keeps the exact same shared
# text + style), then fill B9/C9 with Finland's actual Kapanma-1 dates.
$ws.Range("B9").Copy($ws.Range("D9"))
$ws.Range("C9").Copy($ws.Range("E9"))
$ws.Range("B9").Value = 43906
$ws.Range("C9").Value = 43964

# --- Lithuania (row 16): corrected Kapanma-2 end date
$ws.Range("E16").Value = 44378

# --- Netherlands (row 17): corrected Kapanma-4 end date
$ws.Range("G17").Value = 44542

# --- View/selection state left behind by the editing session
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("G18").Select()
